# Regenerate the handback-status report with a new Handback run:
#   old source "77488849-f5b9-44a8-9f91-6abc60b5c5a0" -> new "de368185-cda9-412a-9fff-1c616eb57438"
#   old source "89e9be6b-9146-4955-aa79-9412f9d98c09" -> new "ffffa84cadc8-772c-4b02-b7be-fe345d7f9392"
# plus refreshed timestamps and xliff-hash filenames.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "77488849-f5b9-44a8-9f91-6abc60b5c5a0"
$oldGuid2 = "89e9be6b-9146-4955-aa79-9412f9d98c09"
$newGuid1 = "de368185-cda9-412a-9fff-1c616eb57438"
$newGuid2 = "ffffa84cadc8-772c-4b02-b7be-fe345d7f9392"

$oldZhXlf1 = "$oldGuid1.624b26f26b67d4c6ff2d736eb3c11cc3da123049.zh-cn.xlf"
$oldZhXlf2 = "$oldGuid2.4487cea072b0a3c5ddbcb88c325b7129c7fc98cf.zh-cn.xlf"
$oldDeXlf1 = "$oldGuid1.624b26f26b67d4c6ff2d736eb3c11cc3da123049.de-de.xlf"
$oldDeXlf2 = "$oldGuid2.4487cea072b0a3c5ddbcb88c325b7129c7fc98cf.de-de.xlf"

$newZhXlf = "$newGuid1.8f4bdc40a08f4e4c42b0ed73a1793094d9db07b7.zh-cn.xlf"
$newDeXlf = "$newGuid1.8f4bdc40a08f4e4c42b0ed73a1793094d9db07b7.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newGuid1.md"
$wsOverview.Range("B2").Value2 = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value2 = "2016-08-24 23:03:41"

$wsOverview.Range("A3").Value2 = "$newGuid2.md"
$wsOverview.Range("B3").Value2 = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value2 = "2016-08-24 23:03:41"

# Rebuild the hyperlinks on B2/B3 so the displayed text matches the refreshed
# file names while keeping the very same external targets.
$ovB2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid1.md"
$ovB3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid2.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovB2Addr, "", "", "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ovB3Addr, "", "", "e2e\$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value2 = "$newGuid1.md"
$wsZh.Range("G2").Value2 = $newZhXlf
$wsZh.Range("H2").Value2 = "2016-08-24 23:03:36"
$wsZh.Range("I2").Value2 = "$newGuid1.md"
$wsZh.Range("J2").Value2 = $newZhXlf
$wsZh.Range("K2").Value2 = "2016-08-24 23:03:53"

$wsZh.Range("A3").Value2 = "$newGuid2.md"
$wsZh.Range("G3").Value2 = $newZhXlf
$wsZh.Range("H3").Value2 = "2016-08-24 23:03:36"
$wsZh.Range("I3").Value2 = "$newGuid2.md"
$wsZh.Range("J3").Value2 = $newZhXlf
$wsZh.Range("K3").Value2 = "2016-08-24 23:03:53"

$zhA2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid1.md"
$zhI2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7d70e96337c3591fc252b8a35f93a0fa73124b71/e2e/$oldGuid1.md"
$zhA3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid2.md"
$zhI3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7d70e96337c3591fc252b8a35f93a0fa73124b71/e2e/$oldGuid2.md"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Addr, "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhI2Addr, "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Addr, "", "", "$newGuid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhI3Addr, "", "", "$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value2 = "$newGuid1.md"
$wsDe.Range("G2").Value2 = $newDeXlf
$wsDe.Range("H2").Value2 = "2016-08-24 23:03:41"
$wsDe.Range("I2").Value2 = "$newGuid1.md"
$wsDe.Range("J2").Value2 = $newDeXlf
$wsDe.Range("K2").Value2 = "2016-08-24 23:04:03"

$wsDe.Range("A3").Value2 = "$newGuid2.md"
$wsDe.Range("G3").Value2 = $newDeXlf
$wsDe.Range("H3").Value2 = "2016-08-24 23:03:41"
$wsDe.Range("I3").Value2 = "$newGuid2.md"
$wsDe.Range("J3").Value2 = $newDeXlf
$wsDe.Range("K3").Value2 = "2016-08-24 23:04:03"

$deA2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid1.md"
$deI2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/23a46f851c1d12ad0d598db188d16d2c83efa757/e2e/$oldGuid1.md"
$deA3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7976737462c1bc379008edfabaa01e2fb0ce5eba/e2e/$oldGuid2.md"
$deI3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/23a46f851c1d12ad0d598db188d16d2c83efa757/e2e/$oldGuid2.md"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Addr, "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deI2Addr, "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Addr, "", "", "$newGuid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deI3Addr, "", "", "$newGuid2.md") | Out-Null

Write-Host "Done."
